$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-orders a block of rows in the "Saldo" export sheet:
#  - A brand-new row for account 004224284 / PRISCILLA / 16030.59 is added right
#    after the ANILSON row (row 2).
#  - The 19 rows that used to sit right after the GILSON row (CESAR .. ROGERIO)
#    are moved up, ending up directly below the new PRISCILLA row.
#  - The GILSON row itself is moved down, to just after the ROGERIO row, and its
#    balance is corrected from 65163.07 to 163.07.
# Net effect: one extra row is inserted, and rows 3-23 end up holding the data
# below (row 2, ANILSON, is untouched; everything from row 24 on is unchanged,
# just shifted down by the single inserted row).

# Insert one blank row at row 3 to make room (this shifts every existing row
# at/after 3 down by one, so nothing below the block is lost).
$ws.Rows.Item(3).Insert()

# Final contents for rows 3 through 23 after the reorder.
$data = @(
    @("004224284", "PRISCILLA", 16030.59),
    @("004207278", "CESAR", 9176.22),
    @("000772433", "MARCELO", 5000),
    @("004313254", "GUSTAVO", 4292),
    @("004368468", "AHMAD", 2766.45),
    @("004213139", "LEONARDO", 2609.78),
    @("004329030", "DANIELA", 940.23),
    @("004392159", "RODRIGO", 900.21),
    @("005696595", "CLUBE", 752.05),
    @("004363260", "LARISSA", 694.83),
    @("005003629", "ANDRE", 650.11),
    @("004855960", "CLERIA", 556.35),
    @("004220849", "DULCE", 503.59),
    @("004432579", "ANA", 446.18),
    @("004508516", "EDUARDO", 364.49),
    @("004355790", "MINEIA", 323.87),
    @("005040864", "ANDRE", 279.96),
    @("003301389", "EDMUNDO", 191.02),
    @("004482090", "CEZAR", 186.91),
    @("004487016", "ROGERIO", 176.96),
    @("004474776", "GILSON", 163.07)
)

# Force column A (the account numbers) to Text so leading zeros survive the
# write, matching the original "inline string" account-number cells.
$ws.Range("A3:A23").NumberFormat = "@"

$row = 3
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}
